$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Determine last used row (data rows 2..51 per known layout, but compute dynamically to be safe)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Merge "First Name" (H) and "Surname" (I) into column H for each data row (rows below header row 1)
for ($r = 2; $r -le $lastRow; $r++) {
    $first = $ws.Cells.Item($r, 8).Value2   # column H
    $last  = $ws.Cells.Item($r, 9).Value2   # column I
    if ($first -eq $null) { $first = "" }
    if ($last -eq $null) { $last = "" }
    $ws.Cells.Item($r, 8).Value2 = "$first $last"
}

# Now delete column I entirely; this shifts J..O left to I..N
$ws.Columns.Item(9).Delete()

# Fix up the header row labels that differ from a pure shift
$ws.Cells.Item(1, 8).Value2 = "Order Name"   # H1
$ws.Cells.Item(1, 11).Value2 = "Mail Box"    # K1 (after shift, this held "Box")

$wb.Save()
